# Extend the Ombudsman appeals table with three more years of data
# (2021, 2022, 2023) in columns R, S, T, mirroring the existing P:Q columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/borders/fonts) of the existing P2:Q5 block onto
# the new R2:S5 block, then copy column Q's formatting onto the new column T.
# PasteSpecial(-4122) == xlPasteFormats: formats only, no values, so cells
# that should stay empty (row 2) remain empty after the paste.
$ws.Range("P2:Q5").Copy()
$ws.Range("R2").PasteSpecial(-4122)

$ws.Range("Q2:Q5").Copy()
$ws.Range("T2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row 3: year headers
$ws.Range("R3").Value = 2021
$ws.Range("S3").Value = 2022
$ws.Range("T3").Value = 2023

# Row 4: "Number of written appeals"
$ws.Range("R4").Value = 4301
$ws.Range("S4").Value = 3690
$ws.Range("T4").Value = 2620

# Row 5: "Number of positively resolved"
$ws.Range("R5").Value = 427
$ws.Range("S5").Value = 280
$ws.Range("T5").Value = 264
